# Update crypto price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.882.40"
$ws.Range("E2").Value = "  -3.83%  "
$ws.Range("D3").Value = "'2.916.78"
$ws.Range("E3").Value = "  -4.26%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'585.92"
$ws.Range("E5").Value = "  -1.70%  "
$ws.Range("D6").Value = "'145.22"
$ws.Range("E6").Value = "  -6.40%  "
$ws.Range("D8").Value = "'0.505"
$ws.Range("E8").Value = "  -2.36%  "
$ws.Range("D9").Value = "'2.916.99"
$ws.Range("E9").Value = "  -4.28%  "
$ws.Range("D10").Value = "'6.83"
$ws.Range("E10").Value = "  -0.03%  "
$ws.Range("D11").Value = "'0.144"
$ws.Range("E11").Value = "  -5.22%  "
$ws.Range("D12").Value = "'0.448"
$ws.Range("E12").Value = "  -4.28%  "
$ws.Range("E13").Value = "  -4.30%  "
$ws.Range("D14").Value = "'33.57"
$ws.Range("E14").Value = "  -6.31%  "
$ws.Range("E15").Value = "  +0.37%  "
$ws.Range("D16").Value = "'3.400.77"
$ws.Range("E16").Value = "  -4.23%  "
$ws.Range("D17").Value = "'60.838.51"
$ws.Range("E17").Value = "  -3.63%  "
$ws.Range("E18").Value = "  -4.65%  "
$ws.Range("D19").Value = "'2.917.92"
$ws.Range("E19").Value = "  -4.13%  "
$ws.Range("D20").Value = "'429.19"
$ws.Range("E20").Value = "  -5.90%  "
$ws.Range("D21").Value = "'13.60"
$ws.Range("E21").Value = "  -5.24%  "
$ws.Range("D22").Value = "'0.683"
$ws.Range("E22").Value = "  -2.66%  "
$ws.Range("D23").Value = "'7.12"
$ws.Range("E23").Value = "  -5.52%  "
$ws.Range("D24").Value = "'80.55"
$ws.Range("E24").Value = "  -3.21%  "
$ws.Range("E25").Value = "  -3.83%  "
$ws.Range("D26").Value = "'10.75"
$ws.Range("E26").Value = "  -5.25%  "
$ws.Range("D27").Value = "'11.97"
$ws.Range("E27").Value = "  -3.60%  "
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("E29").Value = "  +0.34%  "
$ws.Range("D30").Value = "'7.17"
$ws.Range("E30").Value = "  -4.89%  "
$ws.Range("E31").Value = "  -3.37%  "
$ws.Range("E32").Value = "  -3.69%  "
$ws.Range("E33").Value = "  -4.04%  "
$ws.Range("D34").Value = "'0.106"
$ws.Range("E34").Value = "  -4.04%  "
$ws.Range("D35").Value = "'0.0₃0869"
$ws.Range("E35").Value = "  +0.26%  "
$ws.Range("E36").Value = "  -3.45%  "
$ws.Range("E37").Value = "  -5.33%  "
$ws.Range("D38").Value = "'3.01"
$ws.Range("E38").Value = "  -6.91%  "
$ws.Range("E39").Value = "  -3.35%  "
$ws.Range("E40").Value = "  -1.90%  "
$ws.Range("E41").Value = "  -6.11%  "
$ws.Range("D42").Value = "'8.65"
$ws.Range("E42").Value = "  -5.39%  "
$ws.Range("D43").Value = "'0.297"
$ws.Range("E43").Value = "  -2.77%  "
$ws.Range("D44").Value = "'41.08"
$ws.Range("E44").Value = "  -5.83%  "
$ws.Range("D45").Value = "'378.92"
$ws.Range("E45").Value = "  -4.75%  "
$ws.Range("E46").Value = "  -3.56%  "
$ws.Range("D47").Value = "'2.702.64"
$ws.Range("E47").Value = "  -0.97%  "
$ws.Range("D48").Value = "'132.69"
$ws.Range("E48").Value = "  +0.47%  "
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("D50").Value = "'24.65"
$ws.Range("E50").Value = "  +0.29%  "
$ws.Range("D51").Value = "'0.106"
$ws.Range("E51").Value = "  -2.41%  "
